$wb = $excel.ActiveWorkbook

# "ExternalContact" is no longer the active/selected sheet after this edit; reset its
# lingering cell selection (was C9) back to the top-left default before switching away.
$extContactReset = $wb.Worksheets.Item("ExternalContact")
$extContactReset.Activate()
$extContactReset.Range("A1").Select()

# Insert a new worksheet "EngContact" right after "ExternalContact"
$sheets = $wb.Worksheets
$extContact = $sheets.Item("ExternalContact")
$engContact = $sheets.Add([System.Reflection.Missing]::Value, $extContact)
$engContact.Name = "EngContact"

# Populate the new sheet with header + value, matching the other detail sheets' layout
# (assign A2 first so the shared-string table order matches: "Aaron Rosen" before
# "Engagement Contact Name")
$engContact.Range("A2").Value = "Aaron Rosen"
$engContact.Range("A1").Value = "Engagement Contact Name"
$engContact.Range("A1").Font.Bold = $true

# Selection / active state for the new sheet
$engContact.Range("G10").Select()

$wb.Save()
